# Update cryptocurrency list figures (Price / Volume(1h) columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.538.00"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.427.25"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "`'566.05"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "`'145.16"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "`'5.29"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "`'26.80"
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "2.866.26"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "62.315.77"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "2.420.69"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "`'11.21"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "`'6.97"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "`'323.60"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "`'0.999"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "`'67.11"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "`'1.82"
$ws.Range("E24").Value = "  +4.53%  "
$ws.Range("D25").Value = "`'596.49"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "0.0₃0997"
$ws.Range("D28").Value = "2.545.98"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").Value = "`'0.143"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "`'18.72"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").Value = "`'5.35"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "`'147.43"
$ws.Range("E40").Value = "  -3.39%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "`'1.00"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +9.28%  "
$ws.Range("D44").Value = "`'148.32"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "`'20.50"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").Value = "`'0.601"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "`'1.09"
$ws.Range("E51").Value = "  +4.19%  "
